$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting existing rows 47-107 down to 48-108
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new record's values
$ws.Range("A47").Value = 7
$ws.Range("B47").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C47").Value = "Ñuble"
$ws.Range("D47").Value = 44930
$ws.Range("E47").Value = 16
$ws.Range("F47").Value = 100112030
$ws.Range("G47").Value = "Poroto granado"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 41000
$ws.Range("L47").Value = 42000
$ws.Range("M47").Value = 41500
$ws.Range("N47").Value = "$/saco 25 kilos"
$ws.Range("O47").Value = "Región del Maule"
$ws.Range("P47").Value = 1660
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"

Write-Output "done"
